$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = 0.04163557985371009
$ws.Cells.Item(1,2).Value = 0.04163557985371009
$ws.Cells.Item(1,3).Value = 0.03812858876173155
$ws.Cells.Item(1,4).Value = 0.03779284293916666
$ws.Cells.Item(1,5).Value = 0.03812858876173155
$ws.Cells.Item(1,6).Value = 0.03812858876173155
$ws.Cells.Item(1,7).Value = 0.03812858876173155
$ws.Cells.Item(1,8).Value = 0.0355697902479831
$ws.Cells.Item(1,9).Value = 0.0355697902479831
$ws.Cells.Item(1,10).Value = 0.022194926099272
$ws.Cells.Item(1,11).Value = 0.0310160528698903
$ws.Cells.Item(1,12).Value = 0.0310160528698903
$ws.Cells.Item(1,13).Value = 0.03812858876173155
$ws.Cells.Item(1,14).Value = 0.03850847086992486
$ws.Cells.Item(1,15).Value = 0.03850847086992486
$ws.Cells.Item(1,16).Value = 0.0310160528698903
$ws.Cells.Item(1,17).Value = 0.0310160528698903

$ws.Cells.Item(2,1).Value = 0.007401343485758047
$ws.Cells.Item(2,2).Value = 0.007401343485758047
$ws.Cells.Item(2,3).Value = 0.009413718043378096
$ws.Cells.Item(2,4).Value = 0.006718239853483219
$ws.Cells.Item(2,5).Value = 0.009413718043378096
$ws.Cells.Item(2,6).Value = 0.009413718043378096
$ws.Cells.Item(2,7).Value = 0.009413718043378096
$ws.Cells.Item(2,8).Value = 0.006323059178392365
$ws.Cells.Item(2,9).Value = 0.006323059178392365
$ws.Cells.Item(2,10).Value = 0.003945478177052217
$ws.Cells.Item(2,11).Value = 0.007657675933411485
$ws.Cells.Item(2,12).Value = 0.007657675933411485
$ws.Cells.Item(2,13).Value = 0.009413718043378096
$ws.Cells.Item(2,14).Value = 0.006845453360348112
$ws.Cells.Item(2,15).Value = 0.006845453360348112
$ws.Cells.Item(2,16).Value = 0.007657675933411485
$ws.Cells.Item(2,17).Value = 0.007657675933411485

$ws.Cells.Item(3,1).Value = 0.007317156861737104
$ws.Cells.Item(3,2).Value = 0.007317156861737104
$ws.Cells.Item(3,3).Value = 0.01101267627921935
$ws.Cells.Item(3,4).Value = 0.006641823195653999
$ws.Cells.Item(3,5).Value = 0.01101267627921935
$ws.Cells.Item(3,6).Value = 0.01101267627921935
$ws.Cells.Item(3,7).Value = 0.01101267627921935
$ws.Cells.Item(3,8).Value = 0.006251137505423425
$ws.Cells.Item(3,9).Value = 0.006251137505423425
$ws.Cells.Item(3,10).Value = 0.00390060031284912
$ws.Cells.Item(3,11).Value = 0.008958363286135491
$ws.Cells.Item(3,12).Value = 0.008958363286135491
$ws.Cells.Item(3,13).Value = 0.01101267627921935
$ws.Cells.Item(3,14).Value = 0.006767589711753846
$ws.Cells.Item(3,15).Value = 0.006767589711753846
$ws.Cells.Item(3,16).Value = 0.008958363286135491
$ws.Cells.Item(3,17).Value = 0.008958363286135491

$ws.Cells.Item(4,1).Value = 0.007148783613695218
$ws.Cells.Item(4,2).Value = 0.007148783613695218
$ws.Cells.Item(4,3).Value = 0.007176254051528836
$ws.Cells.Item(4,4).Value = 0.006488989879995556
$ws.Cells.Item(4,5).Value = 0.007176254051528836
$ws.Cells.Item(4,6).Value = 0.007176254051528836
$ws.Cells.Item(4,7).Value = 0.007176254051528836
$ws.Cells.Item(4,8).Value = 0.006107294159485542
$ws.Cells.Item(4,9).Value = 0.006107294159485542
$ws.Cells.Item(4,10).Value = 0.003810844584442924
$ws.Cells.Item(4,11).Value = 0.005837590173108594
$ws.Cells.Item(4,12).Value = 0.005837590173108594
$ws.Cells.Item(4,13).Value = 0.007176254051528836
$ws.Cells.Item(4,14).Value = 0.006611862414565312
$ws.Cells.Item(4,15).Value = 0.006611862414565312
$ws.Cells.Item(4,16).Value = 0.005837590173108594
$ws.Cells.Item(4,17).Value = 0.005837590173108594

$ws.Cells.Item(5,1).Value = 0.007148783613695218
$ws.Cells.Item(5,2).Value = 0.007148783613695218
$ws.Cells.Item(5,3).Value = 0.008861747252352734
$ws.Cells.Item(5,4).Value = 0.006488989879995556
$ws.Cells.Item(5,5).Value = 0.008861747252352734
$ws.Cells.Item(5,6).Value = 0.008861747252352734
$ws.Cells.Item(5,7).Value = 0.008861747252352734
$ws.Cells.Item(5,8).Value = 0.006107294159485542
$ws.Cells.Item(5,9).Value = 0.006107294159485542
$ws.Cells.Item(5,10).Value = 0.003810844584442924
$ws.Cells.Item(5,11).Value = 0.007208670192756836
$ws.Cells.Item(5,12).Value = 0.007208670192756836
$ws.Cells.Item(5,13).Value = 0.008861747252352734
$ws.Cells.Item(5,14).Value = 0.006611862414565312
$ws.Cells.Item(5,15).Value = 0.006611862414565312
$ws.Cells.Item(5,16).Value = 0.007208670192756836
$ws.Cells.Item(5,17).Value = 0.007208670192756836

$ws.Cells.Item(6,1).Value = 0.004630180555808607
$ws.Cells.Item(6,2).Value = 0.004630180555808607
$ws.Cells.Item(6,3).Value = 0.005300762804602894
$ws.Cells.Item(6,4).Value = 0.004202840146348175
$ws.Cells.Item(6,5).Value = 0.005300762804602894
$ws.Cells.Item(6,6).Value = 0.005300762804602894
$ws.Cells.Item(6,7).Value = 0.005300762804602894
$ws.Cells.Item(6,8).Value = 0.003955620451524137
$ws.Cells.Item(6,9).Value = 0.003955620451524137
$ws.Cells.Item(6,10).Value = 0.002468237877880805
$ws.Cells.Item(6,11).Value = 0.004311954487109766
$ws.Cells.Item(6,12).Value = 0.004311954487109766
$ws.Cells.Item(6,13).Value = 0.005300762804602894
$ws.Cells.Item(6,14).Value = 0.004282423198675833
$ws.Cells.Item(6,15).Value = 0.004282423198675833
$ws.Cells.Item(6,16).Value = 0.004311954487109766
$ws.Cells.Item(6,17).Value = 0.004311954487109766

$ws.Cells.Item(7,1).Value = 0.004630180555808607
$ws.Cells.Item(7,2).Value = 0.004630180555808607
$ws.Cells.Item(7,3).Value = 0.005211675134796919
$ws.Cells.Item(7,4).Value = 0.004202840146348175
$ws.Cells.Item(7,5).Value = 0.005211675134796919
$ws.Cells.Item(7,6).Value = 0.005211675134796919
$ws.Cells.Item(7,7).Value = 0.005211675134796919
$ws.Cells.Item(7,8).Value = 0.003955620451524137
$ws.Cells.Item(7,9).Value = 0.003955620451524137
$ws.Cells.Item(7,10).Value = 0.002468237877880805
$ws.Cells.Item(7,11).Value = 0.004239485298857755
$ws.Cells.Item(7,12).Value = 0.004239485298857755
$ws.Cells.Item(7,13).Value = 0.005211675134796919
$ws.Cells.Item(7,14).Value = 0.004282423198675833
$ws.Cells.Item(7,15).Value = 0.004282423198675833
$ws.Cells.Item(7,16).Value = 0.004239485298857755
$ws.Cells.Item(7,17).Value = 0.004239485298857755

$ws.Cells.Item(8,1).Value = 0.002315090277904303
$ws.Cells.Item(8,2).Value = 0.002315090277904303
$ws.Cells.Item(8,3).Value = 0.003474450089864613
$ws.Cells.Item(8,4).Value = 0.002101420073174088
$ws.Cells.Item(8,5).Value = 0.003474450089864613
$ws.Cells.Item(8,6).Value = 0.003474450089864613
$ws.Cells.Item(8,7).Value = 0.003474450089864613
$ws.Cells.Item(8,8).Value = 0.001977810225762068
$ws.Cells.Item(8,9).Value = 0.001977810225762068
$ws.Cells.Item(8,10).Value = 0.001234118938940402
$ws.Cells.Item(8,11).Value = 0.002826323532571836
$ws.Cells.Item(8,12).Value = 0.002826323532571836
$ws.Cells.Item(8,13).Value = 0.003474450089864613
$ws.Cells.Item(8,14).Value = 0.002141211599337917
$ws.Cells.Item(8,15).Value = 0.002141211599337917
$ws.Cells.Item(8,16).Value = 0.002826323532571836
$ws.Cells.Item(8,17).Value = 0.002826323532571836

$ws.Cells.Item(9,1).Value = 0.002315090277904303
$ws.Cells.Item(9,2).Value = 0.002315090277904303
$ws.Cells.Item(9,3).Value = 0.001737225044932306
$ws.Cells.Item(9,4).Value = 0.002101420073174088
$ws.Cells.Item(9,5).Value = 0.001737225044932306
$ws.Cells.Item(9,6).Value = 0.001737225044932306
$ws.Cells.Item(9,7).Value = 0.001737225044932306
$ws.Cells.Item(9,8).Value = 0.001977810225762068
$ws.Cells.Item(9,9).Value = 0.001977810225762068
$ws.Cells.Item(9,10).Value = 0.001234118938940402
$ws.Cells.Item(9,11).Value = 0.001413161766285918
$ws.Cells.Item(9,12).Value = 0.001413161766285918
$ws.Cells.Item(9,13).Value = 0.001737225044932306
$ws.Cells.Item(9,14).Value = 0.002141211599337917
$ws.Cells.Item(9,15).Value = 0.002141211599337917
$ws.Cells.Item(9,16).Value = 0.001413161766285918
$ws.Cells.Item(9,17).Value = 0.001413161766285918

$ws.Cells.Item(10,1).Value = 0.06600898496669168
$ws.Cells.Item(10,2).Value = 0.06600898496669168
$ws.Cells.Item(10,3).Value = 0.02505837176621392
$ws.Cells.Item(10,4).Value = 0.05991671570770007
$ws.Cells.Item(10,5).Value = 0.02505837176621392
$ws.Cells.Item(10,6).Value = 0.02505837176621392
$ws.Cells.Item(10,7).Value = 0.02505837176621392
$ws.Cells.Item(10,8).Value = 0.05639229135266317
$ws.Cells.Item(10,9).Value = 0.05639229135266317
$ws.Cells.Item(10,10).Value = 0.03518780207628425
$ws.Cells.Item(10,11).Value = 0.02038396407459801
$ws.Cells.Item(10,12).Value = 0.02038396407459801
$ws.Cells.Item(10,13).Value = 0.02505837176621392
$ws.Cells.Item(10,14).Value = 0.0610512711405567
$ws.Cells.Item(10,15).Value = 0.0610512711405567
$ws.Cells.Item(10,16).Value = 0.02038396407459801
$ws.Cells.Item(10,17).Value = 0.02038396407459801

$wb.Save()
Write-Output "Updated A1:Q10 with realigned k-axis values"
